$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 139; this shifts existing rows 139:207 down to 141:209
$ws.Rows("139:140").Insert()

# Row 139: copy the static (unchanged) descriptive columns from the row below (now row 141,
# which holds what used to be row 139's original data - same market/region/product info)
$ws.Range("A139").Value = 1
$ws.Range("B139").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C139").Value = "Arica y Parinacota"
$ws.Range("D139").Value = 45086
$ws.Range("E139").Value = 15
$ws.Range("F139").Value = "Fruta"
$ws.Range("G139").Value = 100106
$ws.Range("H139").Value = "Oleaginosos"
$ws.Range("I139").Value = 100106002
$ws.Range("J139").Value = "Palta"
$ws.Range("K139").Value = "Hass"
$ws.Range("L139").Value = "Segunda"
$ws.Range("M139").Value = 400
$ws.Range("N139").Value = 24000
$ws.Range("O139").Value = 25000
$ws.Range("P139").Value = 24625
$ws.Range("Q139").Value = "$/bandeja 10 kilos"
$ws.Range("R139").Value = "Perú"
$ws.Range("S139").Value = 2462
$ws.Range("T139").Value = 10

# Row 140
$ws.Range("A140").Value = 1
$ws.Range("B140").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C140").Value = "Arica y Parinacota"
$ws.Range("D140").Value = 45086
$ws.Range("E140").Value = 15
$ws.Range("F140").Value = "Fruta"
$ws.Range("G140").Value = 100106
$ws.Range("H140").Value = "Oleaginosos"
$ws.Range("I140").Value = 100106002
$ws.Range("J140").Value = "Palta"
$ws.Range("K140").Value = "Hass"
$ws.Range("L140").Value = "Segunda"
$ws.Range("M140").Value = 600
$ws.Range("N140").Value = 21000
$ws.Range("O140").Value = 22000
$ws.Range("P140").Value = 21583
$ws.Range("Q140").Value = "$/bandeja 10 kilos"
$ws.Range("R140").Value = "Perú"
$ws.Range("S140").Value = 2158
$ws.Range("T140").Value = 10
